$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.209899999999989
$ws.Range("B21").Value = 5.816099999999994
$ws.Range("B23").Value = 5.623600000000001
$ws.Range("B25").Value = 5.964299999999993
